# Weekly fruit/vegetable price update: insert two new weekly price records
# at the top of the historical data block (rows 192-193), pushing the
# existing historical rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 192:193 - this shifts old rows 192.. down to 194..
$ws.Rows("192:193").Insert()

# Populate the first new row (192) with the new weekly record (Primera)
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(192, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value = "Metropolitana"
$ws.Cells.Item(192, 4).Value = 44726
$ws.Cells.Item(192, 5).Value = 13
$ws.Cells.Item(192, 6).Value = 100112043
$ws.Cells.Item(192, 7).Value = "Pepino ensalada"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 190
$ws.Cells.Item(192, 11).Value = 17000
$ws.Cells.Item(192, 12).Value = 18000
$ws.Cells.Item(192, 13).Value = 17526
$ws.Cells.Item(192, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(192, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(192, 16).Value = 292
$ws.Cells.Item(192, 17).Value = 60
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Populate the second new row (193) with the new weekly record (Segunda)
$ws.Cells.Item(193, 1).Value = 9
$ws.Cells.Item(193, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(193, 3).Value = "Metropolitana"
$ws.Cells.Item(193, 4).Value = 44726
$ws.Cells.Item(193, 5).Value = 13
$ws.Cells.Item(193, 6).Value = 100112043
$ws.Cells.Item(193, 7).Value = "Pepino ensalada"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Segunda"
$ws.Cells.Item(193, 10).Value = 105
$ws.Cells.Item(193, 11).Value = 15000
$ws.Cells.Item(193, 12).Value = 15000
$ws.Cells.Item(193, 13).Value = 15000
$ws.Cells.Item(193, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(193, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(193, 16).Value = 188
$ws.Cells.Item(193, 17).Value = 80
$ws.Cells.Item(193, 18).Value = "Hortaliza"
